# Alteração planilha de tarefas
#
# The task "Revisar documento de requisitos" (old row 15) is removed from the
# list entirely (it was a duplicate-ish leftover blank entry), which shifts
# every subsequent task row up by one. The task that lands on the now-current
# row 15 ("Revisar diagrama de classes") is marked as finished: its Status
# column gets "Pronto" and its Data entrega column gets the delivery date
# 43268 (2018-06-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 15 entirely; rows 16-26 shift up to 15-25.
$ws.Rows("15:15").Delete()

# The new row 15 (previously row 16) needs a Status ("Pronto") and a
# delivery date. Grab the date formatting (numFmtId 14 style) from a sibling
# row (F13) so the new date cell matches the existing "Data entrega" column
# formatting instead of defaulting to General.
$ws.Range("F13").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E15").Value = "Pronto"
$ws.Range("F15").Value = 43268

# Restore a sensible selection/scroll position similar to the source edit.
$ws.Range("E16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
